# Add a "2020" column (I) to the right of the existing data table, carrying
# over the per-row number formatting/border/font from column H and applying
# a new "0.0" number format to the five numeric data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell I4: plain copy of H4's formatting (style already exists,
#     no new number format needed) ---
$ws.Range("H4").Copy($ws.Range("I4"))
$ws.Range("I4").Value = 2020

# --- I5: same formatting as H5 (font 11 row), plus new "0.0" number format ---
$ws.Range("H5").Copy($ws.Range("I5"))
$ws.Range("I5").NumberFormat = "0.0"
$ws.Range("I5").Value = 1.5

# --- I6:I13: same formatting as H6:H13 (font 4 rows), plus "0.0" format ---
$ws.Range("H6").Copy($ws.Range("I6"))
$ws.Range("I6").NumberFormat = "0.0"
$ws.Range("I6").Value = 0.2

$ws.Range("H7").Copy($ws.Range("I7"))
$ws.Range("I7").NumberFormat = "0.0"
$ws.Range("I7").Value = 0.8

$ws.Range("H8").Copy($ws.Range("I8"))
$ws.Range("I8").NumberFormat = "0.0"
$ws.Range("I8").Value = 0.4

$ws.Range("H9").Copy($ws.Range("I9"))
$ws.Range("I9").NumberFormat = "0.0"
$ws.Range("I9").Value = 1.8

$ws.Range("H10").Copy($ws.Range("I10"))
$ws.Range("I10").NumberFormat = "0.0"
$ws.Range("I10").Value = 0.5

$ws.Range("H11").Copy($ws.Range("I11"))
$ws.Range("I11").NumberFormat = "0.0"
$ws.Range("I11").Value = 0.7

$ws.Range("H12").Copy($ws.Range("I12"))
$ws.Range("I12").NumberFormat = "0.0"
$ws.Range("I12").Value = 1.9

$ws.Range("H13").Copy($ws.Range("I13"))
$ws.Range("I13").NumberFormat = "0.0"
$ws.Range("I13").Value = 4.5999999999999996

# --- I14: same border/font as H14 (bottom row), plus "0.0" format, but
#     drop the horizontal/wrap alignment that H14 itself carries so only
#     vertical centering remains (matches the rest of the new column) ---
$ws.Range("H14").Copy($ws.Range("I14"))
$ws.Range("I14").HorizontalAlignment = 1
$ws.Range("I14").WrapText = $false
$ws.Range("I14").NumberFormat = "0.0"
$ws.Range("I14").Value = 0.4

# --- Move the active selection, matching the saved view state ---
$ws.Range("M9").Select()
